$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $textValue) {
    $ws.Range($cellRef).Value = "'" + $textValue
    $ws.Range($cellRef).Style = "Normal"
}

$ws.Range("D2").Value = '29.945.10'
$ws.Range("E2").Value = '  -1.19%  '

$ws.Range("D3").Value = '1.904.34'
$ws.Range("E3").Value = '  -0.22%  '

Set-TextValue 'D4' '1.003'
$ws.Range("E4").Value = '  +0.32%  '

Set-TextValue 'D5' '319.29'
$ws.Range("E5").Value = '  -2.14%  '

Set-TextValue 'D6' '1.002'
$ws.Range("E6").Value = '  +0.19%  '

Set-TextValue 'D7' '0.5029'
$ws.Range("E7").Value = '  -2.69%  '

Set-TextValue 'D8' '0.4035'
$ws.Range("E8").Value = '  +0.07%  '

Set-TextValue 'D9' '0.08244'
$ws.Range("E9").Value = '  -2.81%  '

Set-TextValue 'D10' '41.95'
$ws.Range("E10").Value = '  -1.93%  '

Set-TextValue 'D11' '1.097'
$ws.Range("E11").Value = '  -2.03%  '

Set-TextValue 'D12' '24.14'
$ws.Range("E12").Value = '  +2.17%  '

$ws.Range("D13").Value = '1.908.99'
$ws.Range("E13").Value = '  -0.14%  '

Set-TextValue 'D14' '6.367'
$ws.Range("E14").Value = '  -1.54%  '

Set-TextValue 'D15' '7.191'
$ws.Range("E15").Value = '  -2.23%  '

Set-TextValue 'D16' '1.002'
$ws.Range("E16").Value = '  +0.12%  '

Set-TextValue 'D17' '91.76'
$ws.Range("E17").Value = '  -3.65%  '

Set-TextValue 'D18' '0.00001091'
$ws.Range("E18").Value = '  -2.07%  '

Set-TextValue 'D19' '0.06500'
$ws.Range("E19").Value = '  -2.64%  '

Set-TextValue 'D20' '18.01'
$ws.Range("E20").Value = '  -1.91%  '

Set-TextValue 'D21' '1.002'
$ws.Range("E21").Value = '  +0.16%  '

Set-TextValue 'D22' '5.922'
$ws.Range("E22").Value = '  -1.24%  '

$ws.Range("D23").Value = '29.977.40'
$ws.Range("E23").Value = '  -1.07%  '

Set-TextValue 'D24' '11.20'
$ws.Range("E24").Value = '  -0.83%  '

Set-TextValue 'D25' '2.199'
$ws.Range("E25").Value = '  -1.20%  '

Set-TextValue 'D26' '22.12'
$ws.Range("E26").Value = '  +1.41%  '

$ws.Range("D27").Value = '2.129.12'
$ws.Range("E27").Value = '  +0.03%  '

Set-TextValue 'D28' '161.14'
$ws.Range("E28").Value = '  -0.12%  '

Set-TextValue 'D29' '2.261'
$ws.Range("E29").Value = '  -6.01%  '

Set-TextValue 'D30' '128.46'
$ws.Range("E30").Value = '  -0.89%  '

Set-TextValue 'D31' '1.120'
$ws.Range("E31").Value = '  +1.83%  '

Set-TextValue 'D32' '0.1032'
$ws.Range("E32").Value = '  -2.50%  '

Set-TextValue 'D33' '5.915'
$ws.Range("E33").Value = '  -1.96%  '

Set-TextValue 'D34' '3.798'
$ws.Range("E34").Value = '  +0.66%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D35' '0.02429'
$ws.Range("E35").Value = '  -2.98%  '

$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D36' '5.361'
$ws.Range("E36").Value = '  +1.99%  '

Set-TextValue 'D37' '0.06318'
$ws.Range("E37").Value = '  -3.96%  '

Set-TextValue 'D38' '0.2137'
$ws.Range("E38").Value = '  -3.61%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D39' '1.191'
$ws.Range("E39").Value = '  -3.72%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D40' '0.6466'
$ws.Range("E40").Value = '  -0.90%  '

Set-TextValue 'D41' '8.618'
$ws.Range("E41").Value = '  -2.54%  '

Set-TextValue 'D42' '11.28'
$ws.Range("E42").Value = '  -5.66%  '

Set-TextValue 'D43' '1.202'
$ws.Range("E43").Value = '  -2.82%  '

Set-TextValue 'D44' '2.199'
$ws.Range("E44").Value = '  +6.67%  '

Set-TextValue 'D45' '13.19'
$ws.Range("E45").Value = '  -0.03%  '

Set-TextValue 'D46' '0.5998'
$ws.Range("E46").Value = '  -2.27%  '

Set-TextValue 'D47' '3.629'
$ws.Range("E47").Value = '  -2.27%  '

Set-TextValue 'D48' '122.47'
$ws.Range("E48").Value = '  -2.27%  '

Set-TextValue 'D49' '1.203'
$ws.Range("E49").Value = '  -3.38%  '

Set-TextValue 'D50' '78.22'
$ws.Range("E50").Value = '  -1.55%  '

Set-TextValue 'D51' '1.133'
$ws.Range("E51").Value = '  -2.34%  '
